# Add the four newly-photographed "peeps" rows to the bottom of Table1
# (History of the photo - placed new full size peeps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Columns: Name, No. of peeps, Body type, Peep textures, Head textures,
#          Need work, In garden, Group Photo
$data = @(
    @("Alan",              1, "Full",  $true, $true,  $false, $true,  $false),
    @("John R.",           1, "Head",  $true, $false, $false, $false, $true),
    @("Kris",              1, "Torso", $true, $true,  $false, $true,  $false),
    @("Matt & Lucy w/Lee", 2, "Group", $true, $true,  $false, $true,  $false)
)

foreach ($rowData in $data) {
    $newRow = $tbl.ListRows.Add()
    $r = $newRow.Range.Row

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]
    $ws.Cells.Item($r, 8).Value = $rowData[7]
}

# Update the viewport/selection so H54 is the active cell (matches author's
# final on-screen state after the edit).
$ws.Range("H54").Select()
